$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 9 (Razer), shifting Razer..PcCom down by one
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with Medion (moving it up from its old position)
$ws.Range("A9").Value = "Medion"
$ws.Range("B9").Value = "Medion"

# Remove the old Medion row, now shifted down to row 15, shifting following rows up
$ws.Rows.Item(15).Delete()
